# "delete '-' in cells"
# The workbook used the shared string "-" as a placeholder in a number of
# numeric columns. This edit clears those placeholder cells back to blank
# (keeping their existing cell formatting/style) on every affected sheet,
# then leaves the selection/active sheet the way the author left it when
# they finished editing (VEGFA165_NRP1, with D2:D5 selected).

$wb = $excel.ActiveWorkbook

# --- Adipocyte diameter ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Adipocyte diameter")
$ws1.Range("C4").ClearContents() | Out-Null

# --- CBM thickness ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CBM thickness")
$ws2.Range("D2:E6").ClearContents() | Out-Null
$ws2.Range("D9:E9").ClearContents() | Out-Null

# --- VEGFA165_VEGFR1 ---------------------------------------------------------
$ws3 = $wb.Worksheets.Item("VEGFA165_VEGFR1")
$ws3.Range("D2:D3").ClearContents() | Out-Null

# --- VEGFA165_VEGFR2 ---------------------------------------------------------
$ws4 = $wb.Worksheets.Item("VEGFA165_VEGFR2")
$ws4.Range("D3:D5").ClearContents() | Out-Null

# --- VEGFA165_NRP1 ---------------------------------------------------------
$ws5 = $wb.Worksheets.Item("VEGFA165_NRP1")
$ws5.Range("D2:D5").ClearContents() | Out-Null

# --- Restore each sheet's last selection (as left by the author) -----------
$ws1.Activate() | Out-Null
$ws1.Range("C4").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D9:E9").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("D2:D3").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("D3:D5").Select() | Out-Null

# Final active sheet/tab is VEGFA165_NRP1, selection D2:D5.
$ws5.Activate() | Out-Null
$ws5.Range("D2:D5").Select() | Out-Null
